$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 58) to the FWHM table with results for the
# sg_rr_84_025 2023-12-11 16-27-03 run.
$row = 58

$ws.Cells.Item($row, 1).Value = "sg_rr_84_025 2023-12-11 16-27-03.csv"
$ws.Cells.Item($row, 2).Value = 0.01
$ws.Cells.Item($row, 3).Value = 1000
$ws.Cells.Item($row, 4).Value = 5001
$ws.Cells.Item($row, 5).Value = 1530
$ws.Cells.Item($row, 6).Value = 1570
$ws.Cells.Item($row, 7).Value = 0.5
$ws.Cells.Item($row, 8).Value = "(approx_fsr/2)/wavelength step size"
$ws.Cells.Item($row, 9).Value = 1.7
$ws.Cells.Item($row, 10).Value = 1.1724242424242399
$ws.Cells.Item($row, 11).Value = 0.0059834098769671303
$ws.Cells.Item($row, 12).Value = "yes"
$ws.Cells.Item($row, 13).Value = 0.139397560714696
$ws.Cells.Item($row, 14).Value = 0.0054333187796293501

# Update the view so the newly added row is visible / selected, matching
# the author's saved window state.
$ws.Application.ActiveWindow.ScrollRow = 41
$ws.Range("A59").Select()
